$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet")
$ws.Range("C3").Value = "Sales Order"
